$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set updated forecast-error values (rows 25-52) ---
$ws.Range("J25").Value = -3.899488034340044
$ws.Range("K25").Value = -0.609977258504002

$ws.Range("I26").Value = -3.9
$ws.Range("J26").Value = -0.6

$ws.Range("H27").Value = -3.903239880232202
$ws.Range("I27").Value = -0.615034681294767

$ws.Range("G28").Value = -3.9
$ws.Range("H28").Value = -0.6

$ws.Range("F29").Value = -3.9064445880165
$ws.Range("G29").Value = -0.6210528633565101
$ws.Range("H29").Value = 0.5827860798955309
$ws.Range("I29").Value = 0.5900589890103957
$ws.Range("J29").Value = 0.1909563168408745
$ws.Range("K29").Value = -0.2244894315121023

$ws.Range("E30").Value = -3.905520994775534
$ws.Range("F30").Value = -0.6173747223516166
$ws.Range("G30").Value = 0.5380578365546109
$ws.Range("H30").Value = 0.5922119801925614
$ws.Range("I30").Value = 0.1915127352519465
$ws.Range("J30").Value = -0.2208763097389974

$ws.Range("D31").Value = -3.926584512266816
$ws.Range("E31").Value = -0.3165426857777851
$ws.Range("F31").Value = 0.5434425891838602
$ws.Range("G31").Value = 0.5883235790375505
$ws.Range("H31").Value = 0.1874917024374571
$ws.Range("I31").Value = -0.2246502357367159

$ws.Range("C32").Value = -1.6
$ws.Range("D32").Value = 1.9
$ws.Range("E32").Value = 1.7
$ws.Range("F32").Value = 0.7
$ws.Range("G32").Value = 0.2
$ws.Range("H32").Value = -0.2

$ws.Range("B33").Value = -1.86959775187077
$ws.Range("C33").Value = 0.346995896678686
$ws.Range("D33").Value = 1.233372251821134
$ws.Range("E33").Value = 0.7416915465883819
$ws.Range("F33").Value = 0.3058681882927397
$ws.Range("G33").Value = -0.228687154890392
$ws.Range("H33").Value = 1.36737992645831
$ws.Range("I33").Value = 0.5620836556526756
$ws.Range("J33").Value = 0.09223203051751802
$ws.Range("K33").Value = -0.1080906423023146

$ws.Range("B34").Value = -0.1660949832705663
$ws.Range("C34").Value = 0.3777815514519718
$ws.Range("D34").Value = 0.8724030966903045
$ws.Range("E34").Value = 0.09424305493653341
$ws.Range("F34").Value = -0.2574670601896215
$ws.Range("G34").Value = 1.414665978096479
$ws.Range("H34").Value = 0.5914623413788491
$ws.Range("I34").Value = 0.1151554359168704
$ws.Range("J34").Value = -0.08494058281899805

$ws.Range("B35").Value = 1.797090680270919
$ws.Range("C35").Value = 0.5155367510281608
$ws.Range("D35").Value = -0.4706018608006168
$ws.Range("E35").Value = -0.4700402338977563
$ws.Range("F35").Value = 1.430930696968323
$ws.Range("G35").Value = 0.588985920408935
$ws.Range("H35").Value = 0.07945586353723019
$ws.Range("I35").Value = -0.08793626729908774

$ws.Range("B36").Value = -0.04416977880997519
$ws.Range("C36").Value = -0.3669075164892119
$ws.Range("D36").Value = -0.1984323128299618
$ws.Range("E36").Value = 1.444077178744709
$ws.Range("F36").Value = 0.5756011507456867
$ws.Range("G36").Value = 0.06531475719888355
$ws.Range("H36").Value = -0.1038434665367927

$ws.Range("B37").Value = 0.6901411037136302
$ws.Range("C37").Value = -0.4249506847449998
$ws.Range("D37").Value = 1.654419205737382
$ws.Range("E37").Value = 0.770219776601887
$ws.Range("F37").Value = 0.08688167335395483
$ws.Range("G37").Value = -0.07803537055383716
$ws.Range("H37").Value = -0.6828718050224349
$ws.Range("I37").Value = 0.6147117472547672
$ws.Range("J37").Value = -0.18522021654496
$ws.Range("K37").Value = 0.01485613646721295

$ws.Range("B38").Value = -0.2909814064641099
$ws.Range("C38").Value = 1.605462975383942
$ws.Range("D38").Value = 0.8246063274524341
$ws.Range("E38").Value = 0.3376046120611246
$ws.Range("F38").Value = 0.121816109232675
$ws.Range("G38").Value = -0.5551086423779279
$ws.Range("H38").Value = 0.6113346746433939
$ws.Range("I38").Value = -0.1886353703724037
$ws.Range("J38").Value = 0.01140350262701373

$ws.Range("B39").Value = 1.004708376523435
$ws.Range("C39").Value = 0.7671904148993125
$ws.Range("D39").Value = 0.3654693083104523
$ws.Range("E39").Value = 0.1116880093445911
$ws.Range("F39").Value = -0.5974880634258808
$ws.Range("G39").Value = 0.6095782107951735
$ws.Range("H39").Value = -0.1904109635330261
$ws.Range("I39").Value = 0.0096262989585717

$ws.Range("B40").Value = 0.3619217911506265
$ws.Range("C40").Value = 0.3921640651019387
$ws.Range("D40").Value = 0.425646334868391
$ws.Range("E40").Value = -0.3808891434229855
$ws.Range("F40").Value = 0.7049577299700915
$ws.Range("G40").Value = -0.1979636122625336
$ws.Range("H40").Value = 0.002029828259684563

$ws.Range("B41").Value = 0.07488417233309869
$ws.Range("C41").Value = 0.1460316247142956
$ws.Range("D41").Value = -0.5432921253573036
$ws.Range("E41").Value = 0.5446700041413717
$ws.Range("F41").Value = -0.09386015787386551
$ws.Range("G41").Value = -0.09377001350213054
$ws.Range("H41").Value = -0.4328506734947441
$ws.Range("I41").Value = -0.06647071487539868
$ws.Range("J41").Value = 0.7315413889119743
$ws.Range("K41").Value = 1.029474136650722

$ws.Range("B42").Value = 0.5876823391013496
$ws.Range("C42").Value = -0.6057398398224747
$ws.Range("D42").Value = 0.6542494210714401
$ws.Range("E42").Value = -0.1799629693133837
$ws.Range("F42").Value = -0.02492333770823285
$ws.Range("G42").Value = -0.4177361515802786
$ws.Range("H42").Value = 0.01584405371222172
$ws.Range("I42").Value = 0.8149512355662304
$ws.Range("J42").Value = 1.11399627125113

$ws.Range("B43").Value = -0.3279382706492854
$ws.Range("C43").Value = 0.5581359735930079
$ws.Range("D43").Value = -0.2042459377622568
$ws.Range("E43").Value = -0.01160829470291058
$ws.Range("F43").Value = -0.3719857841332973
$ws.Range("G43").Value = 0.0166722164131696
$ws.Range("H43").Value = 0.7868201425091009
$ws.Range("I43").Value = 1.138126951734013

$ws.Range("B44").Value = 0.4876744685342371
$ws.Range("C44").Value = -0.1521275549781861
$ws.Range("D44").Value = 0.007654823210518802
$ws.Range("E44").Value = -0.2634323168454567
$ws.Range("F44").Value = 0.1335254117779891
$ws.Range("G44").Value = 0.8829196955154583
$ws.Range("H44").Value = 1.240089917752329

$ws.Range("B45").Value = -0.3387895598915543
$ws.Range("C45").Value = -0.02286568139701101
$ws.Range("D45").Value = -0.252137575030208
$ws.Range("E45").Value = 0.1995034391430206
$ws.Range("F45").Value = 0.9436602577277015
$ws.Range("G45").Value = 1.241182925055081
$ws.Range("H45").Value = -0.9450025040107175
$ws.Range("I45").Value = 0.6545759639462473

$ws.Range("B46").Value = 0.1213803088128225
$ws.Range("C46").Value = -0.15027923987009
$ws.Range("D46").Value = 0.1874536666749792
$ws.Range("E46").Value = 0.8870974854954475
$ws.Range("F46").Value = 1.1833997335035
$ws.Range("G46").Value = -0.9487819773954395
$ws.Range("H46").Value = 0.6527813719447992

$ws.Range("B47").Value = -0.3795181086946341
$ws.Range("C47").Value = 0.1094037595619979
$ws.Range("D47").Value = 0.9191112578417829
$ws.Range("E47").Value = 1.181907416094377
$ws.Range("F47").Value = -0.9526883912663854
$ws.Range("G47").Value = 0.6487215976226554

$ws.Range("B48").Value = 0.4863682696630121
$ws.Range("C48").Value = 1.035686727672838
$ws.Range("D48").Value = 1.316280208672642
$ws.Range("E48").Value = -0.8971540607313038
$ws.Range("F48").Value = 0.5811072389458157

$ws.Range("B49").Value = 0.8938767847040198
$ws.Range("C49").Value = 1.180297408122673
$ws.Range("D49").Value = -0.940971934992417
$ws.Range("E49").Value = 0.6685939362799377

$ws.Range("B50").Value = 1.243429160220201
$ws.Range("C50").Value = -0.930588343433577
$ws.Range("D50").Value = 0.672480631954659

$ws.Range("B51").Value = -1.173642127290139
$ws.Range("C51").Value = 0.7060092690864506

$ws.Range("B52").Value = 0.5510376433339623

# --- Remove now-unpopulated trailing cells (rows 45-53) ---
$ws.Range("J45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("H47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("B53").ClearContents()
